# Update "想去人数" (want-to-go count) figures for both the "展览"
# (Exhibition) sheet and the "全部类型" (All types) sheet, matching the
# refreshed scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (rows 2-15) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 4943
$wsExhibit.Range("F4").Value = 340
$wsExhibit.Range("F5").Value = 44
$wsExhibit.Range("F10").Value = 315
$wsExhibit.Range("F11").Value = 252
$wsExhibit.Range("F12").Value = 2967
$wsExhibit.Range("F13").Value = 151
$wsExhibit.Range("F14").Value = 1544
$wsExhibit.Range("F15").Value = 11

# --- Sheet 4: 全部类型 (rows 2-16, offset by one extra row vs sheet 1) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4943
$wsAll.Range("F4").Value = 340
$wsAll.Range("F5").Value = 44
$wsAll.Range("F11").Value = 315
$wsAll.Range("F12").Value = 252
$wsAll.Range("F13").Value = 2967
$wsAll.Range("F14").Value = 151
$wsAll.Range("F15").Value = 1544
$wsAll.Range("F16").Value = 11
